$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 215 (shifts existing rows 215-232 down to 216-233)
$ws.Rows.Item(215).Insert()

# Populate the new row 215 with the new weekly price record.
# Copy the non-changing descriptive columns from the row above (row 214),
# which share the same Mercado/Region/Categoria/etc.
$ws.Range("A215").Value = 7
$ws.Range("B215").Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Range("C215").Value = "Ñuble"
$ws.Range("D215").Value = 45265
$ws.Range("E215").Value = 16
$ws.Range("F215").Value = 100112037
$ws.Range("G215").Value = "Cebollín"
$ws.Range("H215").Value = "Sin especificar"
$ws.Range("I215").Value = "Primera"
$ws.Range("J215").Value = 120
$ws.Range("K215").Value = 5000
$ws.Range("L215").Value = 5000
$ws.Range("M215").Value = 5000
$ws.Range("N215").Value = "$/paquete 36 unidades"
$ws.Range("O215").Value = "Provincia de Diguillín"
$ws.Range("P215").Value = 139
$ws.Range("Q215").Value = 36
$ws.Range("R215").Value = "Hortaliza"
